# Apply updated cryptocurrency price/volume data (scrape refresh).
# Values in column D are stored as text in the source data (e.g. "21.517.20"
# which is not a valid number), so we force text storage via NumberFormat="@"
# around the write, then restore the original cell style so no stray
# formatting is introduced.
function Set-TextValue($Cell, $Value) {
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "21.517.20"
$ws.Range("E2").Value = "  -2.57%  "

Set-TextValue $ws.Range("D3") "1.531.42"
$ws.Range("E3").Value = "  -1.63%  "

Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E5").Value = "  +0.11%  "

Set-TextValue $ws.Range("D6") "288.24"
$ws.Range("E6").Value = "  -1.33%  "

Set-TextValue $ws.Range("D7") "0.3885"
$ws.Range("E7").Value = "  -2.39%  "

Set-TextValue $ws.Range("D8") "0.3175"
$ws.Range("E8").Value = "  -1.69%  "

Set-TextValue $ws.Range("D9") "42.66"
$ws.Range("E9").Value = "  -2.99%  "

Set-TextValue $ws.Range("D10") "0.07151"
$ws.Range("E10").Value = "  -2.23%  "

Set-TextValue $ws.Range("D11") "1.070"
$ws.Range("E11").Value = "  -1.23%  "

Set-TextValue $ws.Range("D13") "5.729"
$ws.Range("E13").Value = "  +0.82%  "

Set-TextValue $ws.Range("D14") "18.17"
$ws.Range("E14").Value = "  -3.91%  "

Set-TextValue $ws.Range("D15") "6.536"
$ws.Range("E15").Value = "  -1.67%  "

Set-TextValue $ws.Range("D16") "1.537.61"
$ws.Range("E16").Value = "  -1.22%  "

Set-TextValue $ws.Range("D17") "0.00001087"
$ws.Range("E17").Value = "  -4.46%  "

Set-TextValue $ws.Range("D18") "0.06608"
$ws.Range("E18").Value = "  +0.15%  "

Set-TextValue $ws.Range("D19") "83.32"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("E20").Value = "  +0.19%  "

Set-TextValue $ws.Range("D21") "6.100"
$ws.Range("E21").Value = "  -3.25%  "

Set-TextValue $ws.Range("D22") "15.40"
$ws.Range("E22").Value = "  -2.05%  "

Set-TextValue $ws.Range("D23") "10.81"
$ws.Range("E23").Value = "  -4.32%  "

Set-TextValue $ws.Range("D24") "2.370"
$ws.Range("E24").Value = "  +0.28%  "

Set-TextValue $ws.Range("D25") "21.496.22"
$ws.Range("E25").Value = "  -2.71%  "

Set-TextValue $ws.Range("D26") "2.373"
$ws.Range("E26").Value = "  -3.02%  "

Set-TextValue $ws.Range("D27") "148.57"
$ws.Range("E27").Value = "  +0.02%  "

Set-TextValue $ws.Range("D28") "18.33"
$ws.Range("E28").Value = "  -1.65%  "

Set-TextValue $ws.Range("D29") "4.835"
$ws.Range("E29").Value = "  -1.24%  "

Set-TextValue $ws.Range("D30") "1.705.94"
$ws.Range("E30").Value = "  -1.45%  "

Set-TextValue $ws.Range("D31") "116.59"
$ws.Range("E31").Value = "  -2.15%  "

Set-TextValue $ws.Range("D32") "6.062"
$ws.Range("E32").Value = "  +4.58%  "

Set-TextValue $ws.Range("D33") "0.9559"
$ws.Range("E33").Value = "  -5.70%  "

Set-TextValue $ws.Range("D34") "0.08009"
$ws.Range("E34").Value = "  -4.09%  "

Set-TextValue $ws.Range("D35") "8.514"
$ws.Range("E35").Value = "  -6.43%  "

Set-TextValue $ws.Range("D36") "5.166"
$ws.Range("E36").Value = "  +0.46%  "

$ws.Range("E37").Value = "  -8.36%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D38") "11.29"
$ws.Range("E38").Value = "  +5.00%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.02197"
$ws.Range("E39").Value = "  -3.32%  "

Set-TextValue $ws.Range("D40") "0.05886"
$ws.Range("E40").Value = "  -3.98%  "

Set-TextValue $ws.Range("D41") "0.2020"
$ws.Range("E41").Value = "  -1.88%  "

Set-TextValue $ws.Range("D42") "1.179"
$ws.Range("E42").Value = "  -3.22%  "

$ws.Range("E43").Value = "  +0.15%  "

Set-TextValue $ws.Range("D44") "0.5742"
$ws.Range("E44").Value = "  -1.87%  "

Set-TextValue $ws.Range("D45") "13.13"
$ws.Range("E45").Value = "  +0.61%  "

Set-TextValue $ws.Range("D46") "3.714"
$ws.Range("E46").Value = "  -1.25%  "

Set-TextValue $ws.Range("D47") "0.5559"
$ws.Range("E47").Value = "  -0.76%  "

Set-TextValue $ws.Range("D48") "1.894"
$ws.Range("E48").Value = "  -1.25%  "

Set-TextValue $ws.Range("D49") "1.160"
$ws.Range("E49").Value = "  +1.69%  "

Set-TextValue $ws.Range("D50") "115.56"
$ws.Range("E50").Value = "  -2.86%  "

Set-TextValue $ws.Range("D51") "0.06678"
$ws.Range("E51").Value = "  -2.49%  "
